$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "(checkbox)" markers in column A for rows 14-16, matching the
# existing pattern used in rows 8-13.
$ws.Range("A14").Value = "(checkbox)"
$ws.Range("A15").Value = "(checkbox)"
$ws.Range("A16").Value = "(checkbox)"

# Update the active selection to D14 (was C12).
$ws.Range("D14").Select()
